# Rerun report after weekend
# Update the numeric result columns (D:J) for rows 2, 3, 5, 7, 9, 10
# with freshly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ D = 1.719909931458051;   E = -0.2039520248260614; F = -0.5852272968399153; G = -0.7854485348414308; H = 0.6762499159617221;  I = -0.5436223594072835; J = -0.2779127231197738 }
    3  = @{ D = 0.09903402696768078; E = 1.244516884631022;   F = 0.3291040908015486;  G = -1.485789282946316;  H = 0.7868775500378195;  I = 0.6466933590286569;  J = -1.620437111431839 }
    5  = @{ D = -2.828694456891071;  E = -3.21069243405172;   F = -11.20394686822662;  G = 12.10898156297891;   H = 4.363515621419242;   I = 4.363512418870462;   J = -3.592691926216328 }
    7  = @{ D = -1.963586432280363;  E = 6.377408859240358;   F = -2.807264789697831;  G = -1.19844024096829;   H = 5.423536811568717;   I = -2.98686097277585;   J = -2.844811869253617 }
    9  = @{ D = -0.350142292383761;  E = 1.453408350316651;   F = 2.443275829108928;   G = 1.188719351958085;   H = 1.823410939967825;   I = 0.5515448522164188;  J = -7.110216070199187 }
    10 = @{ D = -0.5086535384326947; E = -1.046056444018553;  F = -0.4541802921307532; G = 2.035443260118563;   H = 0.4900824914355656;  I = 1.41904260917095;    J = -1.935680723963911 }
}

foreach ($row in $newValues.Keys) {
    $cols = $newValues[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
